$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top to hold per-table (row_count, col_count) summaries
$ws.Rows.Item(1).Insert()

$ws.Range("A1").Value = "(14762, 7)"
$ws.Range("B1").Value = "(772, 6)"
$ws.Range("C1").Value = "(107, 3)"
$ws.Range("D1").Value = "(42, 6)"
$ws.Range("E1").Value = "(9452, 5)"
$ws.Range("F1").Value = "(42, 3)"

$ws.Range("A2").Select()
